$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model detail strings (column B)
$ws.Range("B2").Value = "MLPClassifier(early_stopping=True, learning_rate_init=0.01, max_iter=512,`n              random_state=42)"
$ws.Range("B3").Value = "GaussianNB(var_smoothing=1e-06)"
$ws.Range("B4").Value = "RandomForestClassifier(max_features='sqrt', n_jobs=8, random_state=42)"
$ws.Range("B5").Value = "SVC(cache_size=64, max_iter=1024, random_state=42)"

# Update Mean Acc (column C) and Mean F1 (column E) values
$ws.Range("C2").Value = 0.7480719794344473
$ws.Range("E2").Value = 0.5953078556263269

$ws.Range("C3").Value = 0.2699228791773779
$ws.Range("E3").Value = 0.2385166804521643

$ws.Range("C4").Value = 0.7403598971722365
$ws.Range("E4").Value = 0.576430889312936

$ws.Range("C5").Value = 0.7377892030848329
$ws.Range("E5").Value = 0.4758243500317058

$ws.Range("C6").Value = 0.7429305912596401
$ws.Range("E6").Value = 0.4262536873156342
